$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at row 10 (shifts old rows 10-14 -> 11-15).
#    This makes room for the new "KMR 100" pipe entry and copies the
#    formatting (incl. number format) from the row above, as Excel does.
$ws.Rows("10:10").Insert()

# 2) Break the external workbook link. The formulas in column E referenced
#    an external file ('[1]Rehau Duo'!..., '[1]Rehau Uno'!..., '[1]ISOPLUS KRE einzel'!...);
#    breaking the link converts those formula cells to plain cached values
#    and drops the externalReferences / externalLinks parts entirely.
$wb.BreakLink("/Users/Lars_Goray/Desktop/Waermenetzplanung/Waermenetzplanung/data/Pipe_info.xlsx", 1)

# 3) Header row: column A becomes "Type" (moved from E) and column B keeps
#    "DN", but now the data underneath B is a descriptive text instead of a
#    bare number (see step 4). The header texts themselves do not change.
$ws.Range("A1").Value2 = "Type"
$ws.Range("B1").Value2 = "DN"
$ws.Range("C1").Value2 = "di"
$ws.Range("D1").Value2 = "U-Value"
$ws.Range("E1").Value2 = "U-Value_extra_insulation"
$ws.Range("F1").Value2 = "v_max"
$ws.Range("G1").Value2 = "max_volumeFlow"

# 4) Full data block (rows 2-15). Column A = pipe Type, column B = a text
#    label combining the type family with the DN (e.g. "PEX 20", "KMR 100").
#    Row 10 ("KMR 100") is the newly added pipe; rows 11-15 are the former
#    rows 10-14 shifted down by the insert above.
$rows = @(
  @(2,  "PEX DUO", "PEX 20",  20.399999999999999, 0.129,               0.129,               0.6,                0.19611077980768923),
  @(3,  "PEX DUO", "PEX 25",  26.2,                0.16900000000000001, 0.14299999999999999, 1,                  0.53912871528254436),
  @(4,  "PEX DUO", "PEX 32",  32.6,                0.191,               0.159,               1.1000000000000001, 0.91815872734549941),
  @(5,  "PEX DUO", "PEX 40",  40.799999999999997,  0.17799999999999999, 0.151,               1.2,                1.5688862384615139),
  @(6,  "PEX DUO", "PEX 50",  51.4,                0.21299999999999999, 0.17799999999999999, 1.4,                2.9049867444773314),
  @(7,  "PEX DUO", "PEX 65",  61.4,                0.24299999999999999, 0.24299999999999999, 1.6,                4.737471456130951),
  @(8,  "PEX UNO", "PEX 80",  73.599999999999994,  0.19,                0.16200000000000001, 1.8,                7.6580467833553927),
  @(9,  "PEX UNO", "PEX 100", 90,                  0.27400000000000002, 0.218,               1.9,                12.08727773468673),
  @(10, "KMR UNO", "KMR 100", 107.1,               0.23080000000000001, 0.1943,              1.6,                14.41414231586516),
  @(11, "KMR UNO", "KMR 125", 132.5,               0.26200000000000001, 0.2228,              1.8,                24.819563711063612),
  @(12, "KMR UNO", "KMR 150", 160.30000000000001,  0.30740000000000001, 0.25340000000000001, 2.1,                42.381489974240637),
  @(13, "KMR UNO", "KMR 200", 210.1,               0.33610000000000001, 0.26769999999999999, 2.4,                83.205728598412293),
  @(14, "KMR UNO", "KMR 250", 263,                 0.38229999999999997, 0.29830000000000001, 2.7,                146.67805502290287),
  @(15, "KMR UNO", "KMR 300", 312.7,               0.44169999999999998, 0.34129999999999999, 3,                  230.39173674189982)
)

foreach ($r in $rows) {
  $rowNum = $r[0]
  $ws.Cells.Item($rowNum, 1).Value2 = $r[1]
  $ws.Cells.Item($rowNum, 2).Value2 = $r[2]
  $ws.Cells.Item($rowNum, 3).Value2 = $r[3]
  $ws.Cells.Item($rowNum, 4).Value2 = $r[4]
  $ws.Cells.Item($rowNum, 5).Value2 = $r[5]
  $ws.Cells.Item($rowNum, 6).Value2 = $r[6]
  $ws.Cells.Item($rowNum, 7).Value2 = $r[7]
}

# Note: column B's number format (integer "0" for rows 2-13, "General" for
# the last two rows 14-15) is already inherited correctly from the row
# Insert() above - Excel copies the format of the row above on insert -
# so it is intentionally left untouched here.

# 5) Selection moves to D21 in the saved file.
$ws.Range("D21").Select()
